$d = $word.ActiveDocument

$replacements = @(
    @("2024-02-24 Saturday", "2024-02-25 Sunday"),
    @("518×2=1036", "728×4=2912"),
    @("265×5=1325", "779×9=7011"),
    @("288×5=1440", "405×6=2430"),
    @("668×3=2004", "533×5=2665"),
    @("550×8=4400", "365×4=1460"),
    @("304×6=1824", "422×4=1688"),
    @("669×8=5352", "305×7=2135"),
    @("205×2=410", "531×4=2124"),
    @("733×3=2199", "849×6=5094"),
    @("597×8=4776", "721×7=5047"),
    @("756×3=2268", "837×4=3348"),
    @("747×2=1494", "304×5=1520"),
    @("863×9=7767", "426×9=3834"),
    @("589×8=4712", "448×3=1344"),
    @("833×3=2499", "896×7=6272"),
    @("321×8=2568", "788×3=2364"),
    @("886×2=1772", "331×5=1655"),
    @("135×3=405", "914×4=3656"),
    @("784×6=4704", "152×8=1216"),
    @("561×6=3366", "433×2=866"),
    @("992×7=6944", "460×8=3680"),
    @("820×2=1640", "135×5=675"),
    @("882×4=3528", "643×4=2572"),
    @("150×3=450", "567×8=4536"),
    @("325×2=650", "807×9=7263")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}
